# Updating Old File for core commit 9858844ccecc37046d166185cb936b938f965063
#
# Net effect on the "Rules" sheet (columns A:E = Rule Key, Description, Type,
# Severity, Tags):
#   1. Insert a new row right after row 35 ("BannedPaths") for the rule
#      "CloudServiceIncompatibleWorkflowProcess" at Blocker severity.
#   2. Remove the row "CQRules:CQBP-84--dependencies" (no longer a rule).
#   3. Remove the old "CloudServiceIncompatibleWorkflowProcess" / Major row
#      (it has effectively moved up to Blocker, inserted in step 1).
#   4. Insert a new row right after "PackageOverlaps" or a new rule
#      "IndexDamAssetLucene" at Minor severity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert the new Blocker-severity row just below row 35 ---
$ws.Rows.Item(36).Insert()
$ws.Range("A36").Value = "CloudServiceIncompatibleWorkflowProcess"
$ws.Range("B36").Value = "Usage of Cloud Service Incompatible Workflow Processes"
$ws.Range("C36").Value = "Bug"
$ws.Range("D36").Value = "Blocker"
$ws.Range("E36").Value = "aem,cloud-service-compatibility"

# --- Step 2: delete the "CQRules:CQBP-84--dependencies" row ---
# (row 38 before insert is now row 39)
$ws.Rows.Item(39).Delete()

# --- Step 3: delete the old Major-severity "CloudServiceIncompatibleWorkflowProcess" row ---
# (row 45 before insert/delete is now row 45 again, since steps 1 and 2 cancel out)
$ws.Rows.Item(45).Delete()

# --- Step 4: insert the new Minor-severity "IndexDamAssetLucene" row right after "PackageOverlaps" ---
# (row 70 before all edits is now row 69 ... 35+1-1-1 = 69 net shift of 0 above + insert/delete/delete = -1)
$ws.Rows.Item(70).Insert()
$ws.Range("A70").Value = "IndexDamAssetLucene"
$ws.Range("B70").Value = "Index customizations of the damAssetLucene Oak index should be properly structured."
$ws.Range("C70").Value = "Bug"
$ws.Range("D70").Value = "Minor"
$ws.Range("E70").Value = "aem,cloud-service-compatibility"

# --- Cosmetic: reflect the saved view state (active cell / scroll position) ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 51
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E70").Select() | Out-Null
